# Updated cryptos list on Sat May 18 03:59:37 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to keep text formatting so numeric-looking values
# (e.g. "1.00", "0.999", "581.63") are not auto-converted into real numbers,
# matching the original inline-string cell type.
$ws.Columns("D").NumberFormat = "@"

$ws.Range("D2").Value = "66.978.11"
$ws.Range("E2").Value = "  +2.00%  "

$ws.Range("D3").Value = "3.116.44"
$ws.Range("E3").Value = "  +5.60%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "581.63"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("D6").Value = "173.33"
$ws.Range("E6").Value = "  +6.69%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "3.109.48"
$ws.Range("E8").Value = "  +5.50%  "

$ws.Range("E9").Value = "  +1.46%  "

$ws.Range("D10").Value = "6.51"
$ws.Range("E10").Value = "  -3.21%  "

$ws.Range("E11").Value = "  +3.70%  "

$ws.Range("E12").Value = "  +5.36%  "

$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").Value = "  +1.99%  "

$ws.Range("D14").Value = "37.60"
$ws.Range("E14").Value = "  +8.03%  "

$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("D16").Value = "3.625.66"
$ws.Range("E16").Value = "  +5.33%  "

$ws.Range("D17").Value = "66.992.58"
$ws.Range("E17").Value = "  +2.02%  "

$ws.Range("D18").Value = "7.21"
$ws.Range("E18").Value = "  +2.80%  "

$ws.Range("D19").Value = "3.112.47"
$ws.Range("E19").Value = "  +5.47%  "

$ws.Range("D20").Value = "16.16"
$ws.Range("E20").Value = "  +1.52%  "

$ws.Range("D21").Value = "485.96"
$ws.Range("E21").Value = "  +9.02%  "

$ws.Range("E22").Value = "  +2.93%  "

$ws.Range("D23").Value = "7.55"
$ws.Range("E23").Value = "  +3.42%  "

$ws.Range("D24").Value = "84.23"
$ws.Range("E24").Value = "  +2.51%  "

$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  +7.18%  "

$ws.Range("E26").Value = "  +7.17%  "

$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("D29").Value = "8.03"
$ws.Range("E29").Value = "  -0.73%  "

$ws.Range("E30").Value = "  -4.45%  "

$ws.Range("E31").Value = "  +3.87%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "29.04"
$ws.Range("E32").Value = "  +6.65%  "

$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0000102"
$ws.Range("E33").Value = "  -1.20%  "

$ws.Range("D34").Value = "0.116"
$ws.Range("E34").Value = "  +1.90%  "

$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  +3.48%  "

$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +3.19%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "2.14"
$ws.Range("E38").Value = "  +7.65%  "

$ws.Range("B39").Value = "Arweave"
$ws.Range("C39").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D39").Value = "47.64"
$ws.Range("E39").Value = "  +4.29%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "0.319"
$ws.Range("E40").Value = "  +5.41%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "50.25"
$ws.Range("E41").Value = "  +2.22%  "

$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").Value = "8.70"
$ws.Range("E43").Value = "  +1.79%  "

$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -1.06%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.847.35"
$ws.Range("E45").Value = "  +6.30%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0362"
$ws.Range("E46").Value = "  +2.85%  "

$ws.Range("D47").Value = "383.62"
$ws.Range("E47").Value = "  -0.63%  "

$ws.Range("D48").Value = "134.96"
$ws.Range("E48").Value = "  +1.18%  "

$ws.Range("D50").Value = "25.10"
$ws.Range("E50").Value = "  +5.43%  "

$ws.Range("E51").Value = "  +3.22%  "
